$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-6:
# old serial date 45204 (2023-10-05) -> new serial date 45207 (2023-10-08)
foreach ($r in 2..6) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
